$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G1").Value = "Carry% (>4x)"
$ws.Range("G2").Value = 30
$ws.Range("G3").Value = 30

$s = $wb.Styles.Add("Normal 2")
$s.Font.Name = "Arial"
$s.Font.Size = 11

$ws.Range("G1:G3").Style = "Normal 2"

$ws.Range("G1:G3").Select()
